# Apply the "Add data for 2022-07-22" update:
#  - Rename the sheet tab from "Through 2022-07-13" to "Through 2022-07-14"
#  - Update the running-total column header text in the same way
#  - Bump the August (I8) and Total (I14) figures for the "through" column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-07-14"

$ws.Range("I1").Value = "2022 (through 07-14)"

$ws.Range("I8").Value = 77
$ws.Range("I14").Value = 883
